$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.512729
$ws.Cells.Item(2, 8).Value = 7.538187
$ws.Cells.Item(2, 9).Value = 0.02190726325199687
$ws.Cells.Item(2, 10).Value = 0.02190726325199687
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 98.946724
$ws.Cells.Item(2, 14).Value = 296.840172
$ws.Cells.Item(2, 15).Value = 0.2098009692989996
$ws.Cells.Item(2, 16).Value = 0.2098009692989996
$ws.Cells.Item(2, 17).Value = 248.626302849796
$ws.Cells.Item(2, 18).Value = 2237.636725648164
$ws.Cells.Item(2, 19).Value = 0.004596165064957297
$ws.Cells.Item(2, 20).Value = 0.004596165064957297

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.512729
$ws.Cells.Item(3, 8).Value = 7.538187
$ws.Cells.Item(3, 9).Value = 0.02190726325199687
$ws.Cells.Item(3, 10).Value = 0.02190726325199687
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 163.0062356666667
$ws.Cells.Item(3, 14).Value = 489.018707
$ws.Cells.Item(3, 15).Value = 0.345629090707923
$ws.Cells.Item(3, 16).Value = 0.3456290907079231
$ws.Cells.Item(3, 17).Value = 409.5904955404677
$ws.Cells.Item(3, 18).Value = 3686.314459864209
$ws.Cells.Item(3, 19).Value = 0.007571787477686774
$ws.Cells.Item(3, 20).Value = 0.007571787477686776

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.512729
$ws.Cells.Item(4, 8).Value = 7.538187
$ws.Cells.Item(4, 9).Value = 0.02190726325199687
$ws.Cells.Item(4, 10).Value = 0.02190726325199687
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 65.39610666666668
$ws.Cells.Item(4, 14).Value = 196.18832
$ws.Cells.Item(4, 15).Value = 0.1386621609326595
$ws.Cells.Item(4, 16).Value = 0.1386621609326595
$ws.Cells.Item(4, 17).Value = 164.3226937084267
$ws.Cells.Item(4, 18).Value = 1478.90424337584
$ws.Cells.Item(4, 19).Value = 0.003037708462642526
$ws.Cells.Item(4, 20).Value = 0.003037708462642527

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.512729
$ws.Cells.Item(5, 8).Value = 7.538187
$ws.Cells.Item(5, 9).Value = 0.02190726325199687
$ws.Cells.Item(5, 10).Value = 0.02190726325199687
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 144.2727966666667
$ws.Cells.Item(5, 14).Value = 432.81839
$ws.Cells.Item(5, 15).Value = 0.3059077790604178
$ws.Cells.Item(5, 16).Value = 0.3059077790604179
$ws.Cells.Item(5, 17).Value = 362.5184400954366
$ws.Cells.Item(5, 18).Value = 3262.66596085893
$ws.Cells.Item(5, 19).Value = 0.006701602246710269
$ws.Cells.Item(5, 20).Value = 0.00670160224671027

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 92.89399466666667
$ws.Cells.Item(6, 8).Value = 278.681984
$ws.Cells.Item(6, 9).Value = 0.8098976036382196
$ws.Cells.Item(6, 10).Value = 0.8098976036382197
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 98.946724
$ws.Cells.Item(6, 14).Value = 296.840172
$ws.Cells.Item(6, 15).Value = 0.2098009692989996
$ws.Cells.Item(6, 16).Value = 0.2098009692989996
$ws.Cells.Item(6, 17).Value = 9191.55645154014
$ws.Cells.Item(6, 18).Value = 82724.00806386124
$ws.Cells.Item(6, 19).Value = 0.1699173022762354
$ws.Cells.Item(6, 20).Value = 0.1699173022762355

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 92.89399466666667
$ws.Cells.Item(7, 8).Value = 278.681984
$ws.Cells.Item(7, 9).Value = 0.8098976036382196
$ws.Cells.Item(7, 10).Value = 0.8098976036382197
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 163.0062356666667
$ws.Cells.Item(7, 14).Value = 489.018707
$ws.Cells.Item(7, 15).Value = 0.345629090707923
$ws.Cells.Item(7, 16).Value = 0.3456290907079231
$ws.Cells.Item(7, 17).Value = 15142.30038665274
$ws.Cells.Item(7, 18).Value = 136280.7034798747
$ws.Cells.Item(7, 19).Value = 0.2799241723120037
$ws.Cells.Item(7, 20).Value = 0.2799241723120038

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 92.89399466666667
$ws.Cells.Item(8, 8).Value = 278.681984
$ws.Cells.Item(8, 9).Value = 0.8098976036382196
$ws.Cells.Item(8, 10).Value = 0.8098976036382197
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 65.39610666666668
$ws.Cells.Item(8, 14).Value = 196.18832
$ws.Cells.Item(8, 15).Value = 0.1386621609326595
$ws.Cells.Item(8, 16).Value = 0.1386621609326595
$ws.Cells.Item(8, 17).Value = 6074.905583914099
$ws.Cells.Item(8, 18).Value = 54674.15025522689
$ws.Cells.Item(8, 19).Value = 0.112302151854658
$ws.Cells.Item(8, 20).Value = 0.1123021518546581

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 92.89399466666667
$ws.Cells.Item(9, 8).Value = 278.681984
$ws.Cells.Item(9, 9).Value = 0.8098976036382196
$ws.Cells.Item(9, 10).Value = 0.8098976036382197
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 144.2727966666667
$ws.Cells.Item(9, 14).Value = 432.81839
$ws.Cells.Item(9, 15).Value = 0.3059077790604178
$ws.Cells.Item(9, 16).Value = 0.3059077790604179
$ws.Cells.Item(9, 17).Value = 13402.07640409842
$ws.Cells.Item(9, 18).Value = 120618.6876368858
$ws.Cells.Item(9, 19).Value = 0.2477539771953223
$ws.Cells.Item(9, 20).Value = 0.2477539771953224

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.610639333333333
$ws.Cells.Item(10, 8).Value = 4.831918
$ws.Cells.Item(10, 9).Value = 0.0140423817607685
$ws.Cells.Item(10, 10).Value = 0.0140423817607685
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 98.946724
$ws.Cells.Item(10, 14).Value = 296.840172
$ws.Cells.Item(10, 15).Value = 0.2098009692989996
$ws.Cells.Item(10, 16).Value = 0.2098009692989996
$ws.Cells.Item(10, 17).Value = 159.3674855788773
$ws.Cells.Item(10, 18).Value = 1434.307370209896
$ws.Cells.Item(10, 19).Value = 0.002946105304675824
$ws.Cells.Item(10, 20).Value = 0.002946105304675825

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.610639333333333
$ws.Cells.Item(11, 8).Value = 4.831918
$ws.Cells.Item(11, 9).Value = 0.0140423817607685
$ws.Cells.Item(11, 10).Value = 0.0140423817607685
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 163.0062356666667
$ws.Cells.Item(11, 14).Value = 489.018707
$ws.Cells.Item(11, 15).Value = 0.345629090707923
$ws.Cells.Item(11, 16).Value = 0.3456290907079231
$ws.Cells.Item(11, 17).Value = 262.5442547433362
$ws.Cells.Item(11, 18).Value = 2362.898292690026
$ws.Cells.Item(11, 19).Value = 0.00485345563934794
$ws.Cells.Item(11, 20).Value = 0.004853455639347941

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.610639333333333
$ws.Cells.Item(12, 8).Value = 4.831918
$ws.Cells.Item(12, 9).Value = 0.0140423817607685
$ws.Cells.Item(12, 10).Value = 0.0140423817607685
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 65.39610666666668
$ws.Cells.Item(12, 14).Value = 196.18832
$ws.Cells.Item(12, 15).Value = 0.1386621609326595
$ws.Cells.Item(12, 16).Value = 0.1386621609326595
$ws.Cells.Item(12, 17).Value = 105.3295416441956
$ws.Cells.Item(12, 18).Value = 947.9658747977602
$ws.Cells.Item(12, 19).Value = 0.001947146999589524
$ws.Cells.Item(12, 20).Value = 0.001947146999589524

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.610639333333333
$ws.Cells.Item(13, 8).Value = 4.831918
$ws.Cells.Item(13, 9).Value = 0.0140423817607685
$ws.Cells.Item(13, 10).Value = 0.0140423817607685
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 144.2727966666667
$ws.Cells.Item(13, 14).Value = 432.81839
$ws.Cells.Item(13, 15).Value = 0.3059077790604178
$ws.Cells.Item(13, 16).Value = 0.3059077790604179
$ws.Cells.Item(13, 17).Value = 232.3714410413355
$ws.Cells.Item(13, 18).Value = 2091.34296937202
$ws.Cells.Item(13, 19).Value = 0.004295673817155211
$ws.Cells.Item(13, 20).Value = 0.004295673817155212

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 17.68108066666667
$ws.Cells.Item(14, 8).Value = 53.04324200000001
$ws.Cells.Item(14, 9).Value = 0.154152751349015
$ws.Cells.Item(14, 10).Value = 0.154152751349015
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 98.946724
$ws.Cells.Item(14, 14).Value = 296.840172
$ws.Cells.Item(14, 15).Value = 0.2098009692989996
$ws.Cells.Item(14, 16).Value = 0.2098009692989996
$ws.Cells.Item(14, 17).Value = 1749.485008746403
$ws.Cells.Item(14, 18).Value = 15745.36507871762
$ws.Cells.Item(14, 19).Value = 0.03234139665313101
$ws.Cells.Item(14, 20).Value = 0.03234139665313102

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 17.68108066666667
$ws.Cells.Item(15, 8).Value = 53.04324200000001
$ws.Cells.Item(15, 9).Value = 0.154152751349015
$ws.Cells.Item(15, 10).Value = 0.154152751349015
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 163.0062356666667
$ws.Cells.Item(15, 14).Value = 489.018707
$ws.Cells.Item(15, 15).Value = 0.345629090707923
$ws.Cells.Item(15, 16).Value = 0.3456290907079231
$ws.Cells.Item(15, 17).Value = 2882.126401992011
$ws.Cells.Item(15, 18).Value = 25939.1376179281
$ws.Cells.Item(15, 19).Value = 0.0532796752788846
$ws.Cells.Item(15, 20).Value = 0.05327967527888462

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 17.68108066666667
$ws.Cells.Item(16, 8).Value = 53.04324200000001
$ws.Cells.Item(16, 9).Value = 0.154152751349015
$ws.Cells.Item(16, 10).Value = 0.154152751349015
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 65.39610666666668
$ws.Cells.Item(16, 14).Value = 196.18832
$ws.Cells.Item(16, 15).Value = 0.1386621609326595
$ws.Cells.Item(16, 16).Value = 0.1386621609326595
$ws.Cells.Item(16, 17).Value = 1156.273837259272
$ws.Cells.Item(16, 18).Value = 10406.46453533344
$ws.Cells.Item(16, 19).Value = 0.02137515361576935
$ws.Cells.Item(16, 20).Value = 0.02137515361576936

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 17.68108066666667
$ws.Cells.Item(17, 8).Value = 53.04324200000001
$ws.Cells.Item(17, 9).Value = 0.154152751349015
$ws.Cells.Item(17, 10).Value = 0.154152751349015
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 144.2727966666667
$ws.Cells.Item(17, 14).Value = 432.81839
$ws.Cells.Item(17, 15).Value = 0.3059077790604178
$ws.Cells.Item(17, 16).Value = 0.3059077790604179
$ws.Cells.Item(17, 17).Value = 2550.898955868931
$ws.Cells.Item(17, 18).Value = 22958.09060282038
$ws.Cells.Item(17, 19).Value = 0.04715652580123
$ws.Cells.Item(17, 20).Value = 0.04715652580123001

